$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 (quarter 01-01-2021) with refreshed figures ---
$ws.Range("B75").Value = 396736
$ws.Range("C75").Value = 70080
$ws.Range("D75").Value = 43191
$ws.Range("E75").Value = 26889
$ws.Range("F75").Value = 693
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 689
$ws.Range("I75").Value = 122530
$ws.Range("J75").Value = 97624
$ws.Range("K75").Value = 12698
$ws.Range("L75").Value = 12208
$ws.Range("M75").Value = 158073
$ws.Range("R75").Value = 293695
$ws.Range("S75").Value = 103041
$ws.Range("T75").Value = 98253
$ws.Range("U75").Value = 5402
$ws.Range("V75").Value = 92850
$ws.Range("W75").Value = 4788

# --- Append new row 76 for quarter 01-04-2021 ---
# Column A holds a date-like label that must stay plain text (like the rest
# of column A), so build it as a formula result and convert to a static
# value via copy / paste-special instead of assigning the literal string
# directly (which Excel would otherwise auto-convert to a date serial).
$ws.Range("A76").Formula = "=""01-04-2021"""
$ws.Range("A76").Copy()
$ws.Range("A76").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B76").Value = 396521
$ws.Range("C76").Value = 81044
$ws.Range("D76").Value = 50977
$ws.Range("E76").Value = 30067
$ws.Range("F76").Value = 725
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 721
$ws.Range("I76").Value = 122900
$ws.Range("J76").Value = 98115
$ws.Range("K76").Value = 13207
$ws.Range("L76").Value = 11578
$ws.Range("M76").Value = 146587
$ws.Range("N76").Value = 42682
$ws.Range("O76").Value = 42184
$ws.Range("P76").Value = 498
$ws.Range("Q76").Value = 2583
$ws.Range("R76").Value = 291075
$ws.Range("S76").Value = 105446
$ws.Range("T76").Value = 100614
$ws.Range("U76").Value = 5229
$ws.Range("V76").Value = 95385
$ws.Range("W76").Value = 4832
